$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Num_Inclusions (column C) values to reflect using the 3rd quartile
# instead of the mean for the affected groups.
$ws.Range("C8").Value = 21
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 67
$ws.Range("C13").Value = 17
$ws.Range("C18").Value = 13
$ws.Range("C21").Value = 29
$ws.Range("C25").Value = 0
$ws.Range("C29").Value = 104
$ws.Range("C31").Value = 0
$ws.Range("C37").Value = 0
$ws.Range("C42").Value = 5
$ws.Range("C45").Value = 32
